$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# J12: average of J2:J11 (|S*|/n)
$ws.Range("J12").Formula = "=AVERAGE(J2:J11)"

# A14:A17 labels
$ws.Range("A14").Value = "Average of SW(S*)/SW(OPT)"
$ws.Range("A15").Value = "Average of SC(S*)/SC(OPT)"
$ws.Range("A16").Value = "Worst of SW(S*)/SW(OPT)"
$ws.Range("A17").Value = "Worst of SC(S*)/SC(OPT)"

# B14:B17 formulas
$ws.Range("B14").Formula = "=AVERAGE(N2:N11)"
$ws.Range("B15").Formula = "=AVERAGE(Z2:Z11)"
$ws.Range("B16").Formula = "=MIN(N2:N11)"
$ws.Range("B17").Formula = "=MAX(Z2:Z11)"

# Apply style: bold, size 12, vertical center alignment - build on B14 first,
# then copy formats to B15:B17 so only one new style entry is created.
$b14 = $ws.Range("B14")
$b14.Font.Bold = $true
$b14.Font.Size = 12
$b14.VerticalAlignment = -4108  # xlCenter

$b14.Copy()
$ws.Range("B15:B17").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Row heights for 14-17
$ws.Range("A14:A17").EntireRow.RowHeight = 15.6

# Match the selection left behind by the authoring session
$ws.Range("A14:B17").Select() | Out-Null

# Page setup (paper size 9 = A4, portrait orientation)
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1
